# Auto-applied cell value updates per the authoritative diff.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) gets its changed
# H-N "market price" columns updated to the new scraped values. A few cells
# are newly populated (previously empty) and a few are cleared entirely,
# matching the source diff exactly.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1183.7142
$ws.Range("J32").Value = 1214.3334
$ws.Range("L32").Value = 1214.3334
$ws.Range("N32").Value = -1866.3334
$ws.Range("H43").Value = 1106
$ws.Range("J43").Value = 1157.5
$ws.Range("L43").Value = 1157.5
$ws.Range("N43").Value = -1295.5
$ws.Range("H125").Value = 785005.7
$ws.Range("J125").Value = 1070322.2
$ws.Range("L125").Value = 9632899.799999999
$ws.Range("N125").Value = -9637819.799999999
$ws.Range("H138").Value = 3167.4
$ws.Range("I138").Value = 1597.6154
$ws.Range("J138").Value = 4868
$ws.Range("K138").Value = 4792.8462
$ws.Range("L138").Value = 14604
$ws.Range("M138").Value = 347.1538
$ws.Range("N138").Value = -24884

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2214
$ws.Range("I45").Value = 1025.0667
$ws.Range("K45").Value = 1025.0667
$ws.Range("M45").Value = -648.0667000000001
$ws.Range("H122").Value = 13891713
$ws.Range("I122").Value = 18520428
$ws.Range("J122").Value = 5569
$ws.Range("K122").Value = 55561284
$ws.Range("L122").Value = 16707
$ws.Range("M122").Value = -55558834
$ws.Range("N122").Value = -21607
$ws.Range("H123").Value = 27500
$ws.Range("J123").Value = 27500
$ws.Range("L123").Value = 27500
$ws.Range("N123").Value = -37300
$ws.Range("H135").Value = 29214.5
$ws.Range("J135").Value = 29214.5
$ws.Range("L135").Value = 29214.5
$ws.Range("N135").Value = -39354.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 18030
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 18030
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 18030
$ws.Range("N81").Value = -20152
$ws.Range("H84").Value = 18030
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 18030
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 54090
$ws.Range("N84").Value = -64698
$ws.Range("H132").Value = 44980
$ws.Range("J132").Value = 44980
$ws.Range("L132").Value = 44980
$ws.Range("N132").Value = -55100
$ws.Range("H134").Value = 1716.5264
$ws.Range("I134").Value = 1102.4255
$ws.Range("J134").Value = 4602.8
$ws.Range("K134").Value = 3307.2765
$ws.Range("L134").Value = 13808.4
$ws.Range("M134").Value = -772.2764999999999
$ws.Range("N134").Value = -18878.4
$ws.Range("M81").ClearContents()
$ws.Range("M84").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9613.681
$ws.Range("I31").Value = 1344.2333
$ws.Range("J31").Value = 24206.824
$ws.Range("K31").Value = 1344.2333
$ws.Range("L31").Value = 24206.824
$ws.Range("M31").Value = -1049.2333
$ws.Range("N31").Value = -24796.824
$ws.Range("H34").Value = 9613.681
$ws.Range("I34").Value = 1344.2333
$ws.Range("J34").Value = 24206.824
$ws.Range("K34").Value = 1344.2333
$ws.Range("L34").Value = 24206.824
$ws.Range("M34").Value = -1142.2333
$ws.Range("N34").Value = -24610.824
$ws.Range("H86").Value = 2842.8667
$ws.Range("I86").Value = 3148.6365
$ws.Range("J86").Value = 2002
$ws.Range("K86").Value = 3148.6365
$ws.Range("L86").Value = 2002
$ws.Range("M86").Value = -2025.6365
$ws.Range("N86").Value = -4248
$ws.Range("H89").Value = 2842.8667
$ws.Range("I89").Value = 3148.6365
$ws.Range("J89").Value = 2002
$ws.Range("K89").Value = 15743.1825
$ws.Range("L89").Value = 10010
$ws.Range("M89").Value = -10127.1825
$ws.Range("N89").Value = -21242

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1865.3334
$ws.Range("I5").Value = 874.4
$ws.Range("J5").Value = 2246.4614
$ws.Range("K5").Value = 2623.2
$ws.Range("L5").Value = 6739.3842
$ws.Range("M5").Value = -2511.2
$ws.Range("N5").Value = -6963.3842
$ws.Range("H12").Value = 106.875
$ws.Range("J12").Value = 187.55556
$ws.Range("L12").Value = 562.66668
$ws.Range("N12").Value = -908.66668
$ws.Range("H97").Value = 251
$ws.Range("J97").Value = 218
$ws.Range("L97").Value = 654
$ws.Range("N97").Value = -1646
$ws.Range("H113").Value = 1256.4736
$ws.Range("I113").Value = 3681
$ws.Range("J113").Value = 801.875
$ws.Range("K113").Value = 11043
$ws.Range("L113").Value = 2405.625
$ws.Range("M113").Value = -8873
$ws.Range("N113").Value = -6745.625
$ws.Range("H122").Value = 716.7646999999999
$ws.Range("I122").Value = 534.7143
$ws.Range("K122").Value = 4812.428699999999
$ws.Range("M122").Value = -2362.428699999999
$ws.Range("H132").Value = 1079.9445
$ws.Range("I132").Value = 909.5714
$ws.Range("K132").Value = 8186.1426
$ws.Range("M132").Value = -5656.1426
$ws.Range("H135").Value = 1865.3334
$ws.Range("I135").Value = 874.4
$ws.Range("J135").Value = 2246.4614
$ws.Range("K135").Value = 7869.599999999999
$ws.Range("L135").Value = 20218.1526
$ws.Range("M135").Value = -5334.599999999999
$ws.Range("N135").Value = -25288.1526
$ws.Range("H136").Value = 50509.438
$ws.Range("J136").Value = 131296.17
$ws.Range("L136").Value = 393888.51
$ws.Range("N136").Value = -404088.51

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 218.5
$ws.Range("I2").Value = 41.333332
$ws.Range("J2").Value = 750
$ws.Range("K2").Value = 41.333332
$ws.Range("L2").Value = 750
$ws.Range("M2").Value = 71.666668
$ws.Range("N2").Value = -976
$ws.Range("H70").Value = 6083.2607
$ws.Range("I70").Value = 7217.7334
$ws.Range("J70").Value = 3956.125
$ws.Range("K70").Value = 7217.7334
$ws.Range("L70").Value = 3956.125
$ws.Range("M70").Value = -6947.7334
$ws.Range("N70").Value = -4496.125
$ws.Range("H73").Value = 6083.2607
$ws.Range("I73").Value = 7217.7334
$ws.Range("J73").Value = 3956.125
$ws.Range("K73").Value = 7217.7334
$ws.Range("L73").Value = 3956.125
$ws.Range("M73").Value = -6281.7334
$ws.Range("N73").Value = -5828.125
$ws.Range("H102").Value = 1374.2222
$ws.Range("I102").Value = 1156
$ws.Range("J102").Value = 1548.8
$ws.Range("K102").Value = 1156
$ws.Range("L102").Value = 1548.8
$ws.Range("M102").Value = 466
$ws.Range("N102").Value = -4792.8

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 45274.832
$ws.Range("I40").Value = 57811.11
$ws.Range("J40").Value = 7666
$ws.Range("K40").Value = 57811.11
$ws.Range("L40").Value = 7666
$ws.Range("M40").Value = -57675.11
$ws.Range("N40").Value = -7938
$ws.Range("H46").Value = 27834530
$ws.Range("I46").Value = 62625708
$ws.Range("J46").Value = 1587.8
$ws.Range("K46").Value = 62625708
$ws.Range("L46").Value = 1587.8
$ws.Range("M46").Value = -62625520
$ws.Range("N46").Value = -1963.8
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("H122").Value = 2970.8948
$ws.Range("I122").Value = 1800.6666
$ws.Range("J122").Value = 3511
$ws.Range("K122").Value = 5401.9998
$ws.Range("L122").Value = 10533
$ws.Range("M122").Value = -2951.9998
$ws.Range("N122").Value = -15433
$ws.Range("H136").Value = 3116.2856
$ws.Range("J136").Value = 6679
$ws.Range("L136").Value = 20037
$ws.Range("N136").Value = -25137
$ws.Range("N64").ClearContents()
$ws.Range("N67").ClearContents()

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 21000
$ws.Range("J63").Value = 21000
$ws.Range("L63").Value = 21000
$ws.Range("N63").Value = -22248
$ws.Range("H66").Value = 21000
$ws.Range("J66").Value = 21000
$ws.Range("L66").Value = 63000
$ws.Range("N66").Value = -69240
$ws.Range("H122").Value = 33883.227
$ws.Range("I122").Value = 40003.08
$ws.Range("J122").Value = 2060
$ws.Range("K122").Value = 120009.24
$ws.Range("L122").Value = 6180
$ws.Range("M122").Value = -117559.24
$ws.Range("N122").Value = -11080
